# Regenerate merged AHB files
# 1. Rename the "_old" / "_new" header-name suffixes to "_FV2304" / "_FV2310"
# 2. Wrap the data range in a native Excel Table ("Table1")
# 3. Freeze the header row (row 1) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the existing A1:U56 range into a real ListObject/Table so it matches
# the regenerated-report layout (autofilter + table part).
$dataRange = $ws.Range("A1:U56")
$lo = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$lo.Name = "Table1"

# Freeze the header row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
